$d = $word.ActiveDocument

$replacements = @(
    @{old="960×7="; new="564×7="},
    @{old="765×5="; new="379×3="},
    @{old="409×4="; new="298×3="},
    @{old="364×2="; new="403×6="},
    @{old="817×2="; new="481×6="},
    @{old="574×6="; new="527×2="},
    @{old="958×9="; new="740×4="},
    @{old="698×2="; new="849×8="},
    @{old="579×2="; new="426×5="},
    @{old="985×8="; new="111×6="},
    @{old="936×5="; new="723×3="},
    @{old="445×4="; new="375×2="},
    @{old="995×5="; new="966×2="},
    @{old="491×3="; new="358×9="},
    @{old="176×9="; new="575×5="},
    @{old="440×2="; new="987×2="},
    @{old="807×3="; new="175×4="},
    @{old="304×6="; new="583×6="},
    @{old="111×7="; new="488×9="},
    @{old="670×8="; new="429×6="},
    @{old="473×9="; new="901×9="},
    @{old="589×7="; new="293×5="},
    @{old="589×6="; new="707×7="},
    @{old="494×3="; new="352×6="},
    @{old="711×3="; new="710×3="}
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
